# Apply scheduled market-data refresh to the Leve profit sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the rows whose underlying market data changed since the last run.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 113 (Leve Item ID G113=27775)
$ws.Range("H113").Value = 2566.6667
$ws.Range("I113").Value = 2514.2856
$ws.Range("J113").Value = 2750
$ws.Range("K113").Value = 2514.2856
$ws.Range("L113").Value = 2750
$ws.Range("M113").Value = 739.7143999999998
$ws.Range("N113").Value = -9258
# Row 138 (Leve Item ID G138=44169)
$ws.Range("H138").Value = 2347.1094
$ws.Range("I138").Value = 2234.5
$ws.Range("J138").Value = 2391.1738
$ws.Range("K138").Value = 6703.5
$ws.Range("L138").Value = 7173.5214
$ws.Range("M138").Value = -1563.5
$ws.Range("N138").Value = -17453.5214

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45 (Leve Item ID G45=27714)
$ws.Range("H45").Value = 4088.6667
$ws.Range("I45").Value = 3756
$ws.Range("J45").Value = 4255
$ws.Range("K45").Value = 3756
$ws.Range("L45").Value = 4255
$ws.Range("M45").Value = -3379
$ws.Range("N45").Value = -5009
# Row 61 (Leve Item ID G61=43999)
$ws.Range("H61").Value = 10103182
$ws.Range("I61").Value = 19609164
$ws.Range("J61").Value = 3076.75
$ws.Range("K61").Value = 19609164
$ws.Range("L61").Value = 3076.75
$ws.Range("M61").Value = -19608952
$ws.Range("N61").Value = -3500.75
# Row 132 (Leve Item ID G132=43997)
$ws.Range("H132").Value = 4362.838
$ws.Range("I132").Value = 4774.115
$ws.Range("J132").Value = 3390.7273
$ws.Range("K132").Value = 14322.345
$ws.Range("L132").Value = 10172.1819
$ws.Range("M132").Value = -11792.345
$ws.Range("N132").Value = -15232.1819
# Row 136 (Leve Item ID G136=43999)
$ws.Range("H136").Value = 10103182
$ws.Range("I136").Value = 19609164
$ws.Range("J136").Value = 3076.75
$ws.Range("K136").Value = 58827492
$ws.Range("L136").Value = 9230.25
$ws.Range("M136").Value = -58824942
$ws.Range("N136").Value = -14330.25

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID G16=27691)
$ws.Range("H16").Value = 556
$ws.Range("I16").Value = 496.66666
$ws.Range("J16").Value = 645
$ws.Range("K16").Value = 496.66666
$ws.Range("L16").Value = 645
$ws.Range("M16").Value = -209.66666
$ws.Range("N16").Value = -1219
# Row 31 (Leve Item ID G31=44023)
$ws.Range("H31").Value = 3168.4666
$ws.Range("I31").Value = 1511.0625
$ws.Range("J31").Value = 3771.1592
$ws.Range("K31").Value = 1511.0625
$ws.Range("L31").Value = 3771.1592
$ws.Range("M31").Value = -1216.0625
$ws.Range("N31").Value = -4361.1592
# Row 34 (Leve Item ID G34=44023)
$ws.Range("H34").Value = 3168.4666
$ws.Range("I34").Value = 1511.0625
$ws.Range("J34").Value = 3771.1592
$ws.Range("K34").Value = 1511.0625
$ws.Range("L34").Value = 3771.1592
$ws.Range("M34").Value = -1309.0625
$ws.Range("N34").Value = -4175.1592
# Row 58 (Leve Item ID G58=44021)
$ws.Range("H58").Value = 3506.6667
$ws.Range("I58").Value = 3697.5
$ws.Range("J58").Value = 3125
$ws.Range("K58").Value = 3697.5
$ws.Range("L58").Value = 3125
$ws.Range("M58").Value = -3494.5
$ws.Range("N58").Value = -3531
# Row 94 (Leve Item ID G94=32934)
$ws.Range("H94").Value = 1836.7142
$ws.Range("J94").Value = 1836.7142
$ws.Range("L94").Value = 1836.7142
$ws.Range("N94").Value = -2738.7142
# Row 99 (Leve Item ID G99=36198)
$ws.Range("H99").Value = 1772.88
$ws.Range("I99").Value = 1153.6666
$ws.Range("J99").Value = 1968.421
$ws.Range("K99").Value = 1153.6666
$ws.Range("L99").Value = 1968.421
$ws.Range("M99").Value = 344.3334
$ws.Range("N99").Value = -4964.421
# Row 105 (Leve Item ID G105=19928)
$ws.Range("H105").Value = 378.3
$ws.Range("I105").Value = 349.125
$ws.Range("J105").Value = 495
$ws.Range("K105").Value = 349.125
$ws.Range("L105").Value = 495
$ws.Range("M105").Value = 1397.875
$ws.Range("N105").Value = -3989
# Row 107 (Leve Item ID G107=27689)
$ws.Range("H107").Value = 2083798.5
$ws.Range("I107").Value = 3906575.8
$ws.Range("J107").Value = 624.4286
$ws.Range("K107").Value = 3906575.8
$ws.Range("L107").Value = 624.4286
$ws.Range("M107").Value = -3904655.8
$ws.Range("N107").Value = -4464.4286
# Row 113 (Leve Item ID G113=27691)
$ws.Range("H113").Value = 556
$ws.Range("I113").Value = 496.66666
$ws.Range("J113").Value = 645
$ws.Range("K113").Value = 496.66666
$ws.Range("L113").Value = 645
$ws.Range("M113").Value = 1673.33334
$ws.Range("N113").Value = -4985
# Row 126 (Leve Item ID G126=36198)
$ws.Range("H126").Value = 1772.88
$ws.Range("I126").Value = 1153.6666
$ws.Range("J126").Value = 1968.421
$ws.Range("K126").Value = 3460.9998
$ws.Range("L126").Value = 5905.263
$ws.Range("M126").Value = -990.9998000000001
$ws.Range("N126").Value = -10845.263
# Row 136 (Leve Item ID G136=44021)
$ws.Range("H136").Value = 3506.6667
$ws.Range("I136").Value = 3697.5
$ws.Range("J136").Value = 3125
$ws.Range("K136").Value = 11092.5
$ws.Range("L136").Value = 9375
$ws.Range("M136").Value = -8542.5
$ws.Range("N136").Value = -14475

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68 (Leve Item ID G68=12895)
$ws.Range("H68").Value = 1414.2872
$ws.Range("I68").Value = 791.2593000000001
$ws.Range("K68").Value = 2373.7779
$ws.Range("M68").Value = -1562.7779
# Row 71 (Leve Item ID G71=12895)
$ws.Range("H71").Value = 1414.2872
$ws.Range("I71").Value = 791.2593000000001
$ws.Range("K71").Value = 7121.3337
$ws.Range("M71").Value = -3065.3337
# Row 112 (Leve Item ID G112=27855)
$ws.Range("H112").Value = 6750
$ws.Range("I112").Value = 8000
$ws.Range("J112").Value = 6666.6665
$ws.Range("K112").Value = 24000
$ws.Range("L112").Value = 19999.9995
$ws.Range("M112").Value = -22892
$ws.Range("N112").Value = -22215.9995
# Row 122 (Leve Item ID G122=36078)
$ws.Range("H122").Value = 5508.35
$ws.Range("J122").Value = 15027.286
$ws.Range("L122").Value = 135245.574
$ws.Range("N122").Value = -140145.574
# Row 131 (Leve Item ID G131=36060)
$ws.Range("H131").Value = 1180
$ws.Range("J131").Value = 1192.5
$ws.Range("L131").Value = 3577.5
$ws.Range("N131").Value = -13657.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102 (Leve Item ID G102=36169)
$ws.Range("H102").Value = 2102.5
$ws.Range("I102").Value = 2004.375
$ws.Range("K102").Value = 2004.375
$ws.Range("M102").Value = -382.375
# Row 113 (Leve Item ID G113=27710)
$ws.Range("H113").Value = 1668.8182
$ws.Range("I113").Value = 1466.9166
$ws.Range("J113").Value = 1911.1
$ws.Range("K113").Value = 1466.9166
$ws.Range("L113").Value = 1911.1
$ws.Range("M113").Value = 703.0834
$ws.Range("N113").Value = -6251.1
# Row 121 (Leve Item ID G121=26338)
$ws.Range("H121").Value = 39800
$ws.Range("J121").Value = 39800
$ws.Range("L121").Value = 39800
$ws.Range("N121").Value = -43294
# Row 122 (Leve Item ID G122=36182)
$ws.Range("H122").Value = 5740.5835
$ws.Range("J122").Value = 5388
$ws.Range("L122").Value = 16164
$ws.Range("N122").Value = -21064

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID G7=36249)
$ws.Range("H7").Value = 500001500
$ws.Range("I7").Value = 500001500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 500001500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -500001388
$ws.Range("N7").ClearContents()
# Row 40 (Leve Item ID G40=36248)
$ws.Range("H40").Value = 333336670
$ws.Range("I40").Value = 1000000000
$ws.Range("K40").Value = 1000000000
$ws.Range("M40").Value = -999999864
# Row 55 (Leve Item ID G55=5284)
$ws.Range("H55").Value = 395.35483
$ws.Range("I55").Value = 138.44444
$ws.Range("J55").Value = 751.0769
$ws.Range("K55").Value = 138.44444
$ws.Range("L55").Value = 751.0769
$ws.Range("M55").Value = 34.55556000000001
$ws.Range("N55").Value = -1097.0769
# Row 122 (Leve Item ID G122=36247)
$ws.Range("H122").Value = 5050
$ws.Range("I122").Value = 2300
$ws.Range("J122").Value = 5966.6665
$ws.Range("K122").Value = 6900
$ws.Range("L122").Value = 17899.9995
$ws.Range("M122").Value = -4450
$ws.Range("N122").Value = -22799.9995
# Row 126 (Leve Item ID G126=36249)
$ws.Range("H126").Value = 500001500
$ws.Range("I126").Value = 500001500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 1500004500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1500002030
$ws.Range("N126").ClearContents()
# Row 132 (Leve Item ID G132=44058)
$ws.Range("H132").Value = 2961.3901
$ws.Range("I132").Value = 2953.5264
$ws.Range("J132").Value = 2968.182
$ws.Range("K132").Value = 8860.5792
$ws.Range("L132").Value = 8904.545999999998
$ws.Range("M132").Value = -6330.5792
$ws.Range("N132").Value = -13964.546

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107 (Leve Item ID G107=27746)
$ws.Range("H107").Value = 548.6087
$ws.Range("I107").Value = 416.9
$ws.Range("J107").Value = 649.9231
$ws.Range("K107").Value = 1250.7
$ws.Range("L107").Value = 1949.7693
$ws.Range("M107").Value = 669.3000000000002
$ws.Range("N107").Value = -5789.7693
# Row 122 (Leve Item ID G122=36208)
$ws.Range("H122").Value = 2426
$ws.Range("I122").Value = 2069.3333
$ws.Range("K122").Value = 6207.999899999999
$ws.Range("M122").Value = -3757.999899999999
# Row 132 (Leve Item ID G132=44029)
$ws.Range("H132").Value = 5504967
$ws.Range("I132").Value = 1702.4073
$ws.Range("J132").Value = 11219896
$ws.Range("K132").Value = 5107.2219
$ws.Range("L132").Value = 33659688
$ws.Range("M132").Value = -2577.2219
$ws.Range("N132").Value = -33664748
# Row 136 (Leve Item ID G136=44031)
$ws.Range("H136").Value = 2993.7046
$ws.Range("I136").Value = 2673.5757
$ws.Range("J136").Value = 3954.0908
$ws.Range("K136").Value = 8020.7271
$ws.Range("L136").Value = 11862.2724
$ws.Range("M136").Value = -5470.7271
$ws.Range("N136").Value = -16962.2724
